$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells: _old -> _FV2410, _new -> _FV2504 (column K "diff" unchanged)
$oldHeaders = @(
  "Segmentname_FV2410",
  "Segmentgruppe_FV2410",
  "Segment_FV2410",
  "Datenelement_FV2410",
  "Segment ID_FV2410",
  "Code_FV2410",
  "Qualifier_FV2410",
  "Beschreibung_FV2410",
  "Bedingungsausdruck_FV2410",
  "Bedingung_FV2410"
)

$newHeaders = @(
  "Segmentname_FV2504",
  "Segmentgruppe_FV2504",
  "Segment_FV2504",
  "Datenelement_FV2504",
  "Segment ID_FV2504",
  "Code_FV2504",
  "Qualifier_FV2504",
  "Beschreibung_FV2504",
  "Bedingungsausdruck_FV2504",
  "Bedingung_FV2504"
)

for ($i = 0; $i -lt 10; $i++) {
  $ws.Cells.Item(1, $i + 1).Value = $oldHeaders[$i]
}

for ($i = 0; $i -lt 10; $i++) {
  $ws.Cells.Item(1, $i + 12).Value = $newHeaders[$i]
}

# Freeze top row
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Add table covering A1:U55
$rng = $ws.Range("A1:U55")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $rng, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

$wb.Save()
